$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'60.566.95"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.02%  "
$ws.Range("D3").Value = "'2.386.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -4.38%  "
$ws.Range("E4").Value = "  +0.21%  "
$ws.Range("D5").Value = "'536.20"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.58%  "
$ws.Range("D6").Value = "'139.20"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -5.15%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "'0.570"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -6.68%  "
$ws.Range("D9").Value = "'2.386.48"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -4.38%  "
$ws.Range("E10").Value = "  -3.89%  "
$ws.Range("E11").Value = "  -0.14%  "
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "'0.337"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.95%  "
$ws.Range("D14").Value = "'25.14"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.58%  "
$ws.Range("D15").Value = "'2.825.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -4.02%  "
$ws.Range("D16").Value = "'60.955.13"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.23%  "
$ws.Range("E17").Value = "  -0.99%  "
$ws.Range("D18").Value = "'2.387.12"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.21%  "
$ws.Range("D19").Value = "'10.50"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.39%  "
$ws.Range("D20").Value = "'6.68"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -4.69%  "
$ws.Range("D21").Value = "'4.03"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -3.46%  "
$ws.Range("D22").Value = "'311.03"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.90%  "
$ws.Range("D23").Value = "'1.00"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("D24").Value = "'1.76"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.93%  "
$ws.Range("D25").Value = "'62.47"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.10%  "
$ws.Range("E26").Value = "  +0.14%  "
$ws.Range("D27").Value = "'2.506.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -4.50%  "
$ws.Range("D28").Value = "'0.0₃0899"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -11.33%  "
$ws.Range("D29").Value = "'7.51"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.16%  "
$ws.Range("E30").Value = "  -5.60%  "
$ws.Range("D31").Value = "'7.92"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.16%  "
$ws.Range("D32").Value = "'500.50"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -6.38%  "
$ws.Range("E33").Value = "  -4.65%  "
$ws.Range("D34").Value = "'1.81"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -4.24%  "
$ws.Range("E35").Value = "  -2.37%  "
$ws.Range("D36").Value = "'1.00"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.45%  "
$ws.Range("D37").Value = "'5.41"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -7.98%  "
$ws.Range("D38").Value = "'4.56"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.20%  "
$ws.Range("D39").Value = "'0.369"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.21%  "
$ws.Range("D40").Value = "'17.89"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.03%  "
$ws.Range("D41").Value = "'1.00"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.07%  "
$ws.Range("D42").Value = "'136.75"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.88%  "
$ws.Range("E43").Value = "  -1.18%  "
$ws.Range("D44").Value = "'40.22"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.40%  "
$ws.Range("B45").Value = "Aave"
$ws.Range("C45").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D45").Value = "'139.15"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -6.49%  "
$ws.Range("B46").Value = "dogwifhat"
$ws.Range("C46").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D46").Value = "'2.11"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -8.29%  "
$ws.Range("D47").Value = "'3.48"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.54%  "
$ws.Range("D48").Value = "'19.88"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.08%  "
$ws.Range("E49").Value = "  -4.59%  "
$ws.Range("E50").Value = "  -2.82%  "
$ws.Range("D51").Value = "'0.0916"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.89%  "
